$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 3.85
$ws.Range("F5").Value = 2.34
$ws.Range("G5").Value = 2.66
$ws.Range("I5").Value = 3
$ws.Range("K5").Value = 4.5
$ws.Range("T5").Value = 1.54
$ws.Range("F6").Value = 2.58
$ws.Range("H6").Value = 2.66
$ws.Range("I6").Value = 3.5
$ws.Range("J6").Value = 2.48
$ws.Range("K6").Value = 4.3
$ws.Range("Q6").Value = 2.06
$ws.Range("N7").Value = 4.7
$ws.Range("P7").Value = 2.28
$ws.Range("T7").Value = 1.63
$ws.Range("U7").Value = 2.3
$ws.Range("Z7").Value = 1000
$ws.Range("AF7").Value = 16.5
$ws.Range("AG7").Value = 12.5
$ws.Range("AH7").Value = 21
$ws.Range("AN7").Value = 11.5
$ws.Range("I8").Value = 3.85
$ws.Range("F13").Value = 1.91
$ws.Range("G13").Value = 1.97
$ws.Range("P13").Value = 1.7
$ws.Range("Q13").Value = 2.2
$ws.Range("G15").Value = 2.58
$ws.Range("G17").Value = 2.72
$ws.Range("H17").Value = 2.9
$ws.Range("I17").Value = 3.65
$ws.Range("J17").Value = 3.3
$ws.Range("K17").Value = 5.6
$ws.Range("P17").Value = 2.04
$ws.Range("Q17").Value = 1.69
$ws.Range("V17").Value = 1.37
$ws.Range("W17").Value = 1.58
$ws.Range("G19").Value = 1.74
$ws.Range("R19").Value = 1.64
$ws.Range("S19").Value = 2.36
$ws.Range("W19").Value = 2.36
$ws.Range("AC19").Value = 11.5
$ws.Range("AO19").Value = 50
$ws.Range("G20").Value = 2.9
$ws.Range("H20").Value = 2.9
$ws.Range("K20").Value = 3.4
$ws.Range("V20").Value = 1.48
$ws.Range("AE20").Value = 36
$ws.Range("AI20").Value = 980
$ws.Range("H22").Value = 2.1
$ws.Range("G24").Value = 2.42
$ws.Range("H24").Value = 3.1
$ws.Range("K24").Value = 4.1
$ws.Range("Q24").Value = 1.56
$ws.Range("Q25").Value = 2
$ws.Range("F26").Value = 1.44
$ws.Range("H29").Value = 4.1
$ws.Range("K29").Value = 6.2
$ws.Range("G30").Value = 2.9
$ws.Range("H30").Value = 2.42
$ws.Range("K30").Value = 4.5
$ws.Range("P30").Value = 2.64
$ws.Range("N31").Value = 1.1
$ws.Range("R31").Value = 1.63
$ws.Range("G32").Value = 9
$ws.Range("F37").Value = 2.08
$ws.Range("AM37").Value = 160
$ws.Range("Q38").Value = 1.54
$ws.Range("G39").Value = 3.4
$ws.Range("H39").Value = 2.34
$ws.Range("I39").Value = 2.62
$ws.Range("M41").Value = 1.09
$ws.Range("N41").Value = 2.92
$ws.Range("O41").Value = 1.41
$ws.Range("Q41").Value = 2.22
$ws.Range("S41").Value = 4.2
$ws.Range("T41").Value = 1.87
$ws.Range("U41").Value = 1.92
$ws.Range("V41").Value = 1.49
$ws.Range("X41").Value = 11
$ws.Range("Y41").Value = 10
$ws.Range("Z41").Value = 19
$ws.Range("AA41").Value = 50
$ws.Range("AB41").Value = 10.5
$ws.Range("AC41").Value = 7.6
$ws.Range("AD41").Value = 13.5
$ws.Range("AE41").Value = 38
$ws.Range("AF41").Value = 19.5
$ws.Range("AG41").Value = 14
$ws.Range("AH41").Value = 980
$ws.Range("AI41").Value = 55
$ws.Range("AJ41").Value = 55
$ws.Range("AK41").Value = 40
$ws.Range("AL41").Value = 60
$ws.Range("AN41").Value = 42
$ws.Range("AO41").Value = 40
$ws.Range("F42").Value = 1.99
$ws.Range("G42").Value = 2.18
$ws.Range("H42").Value = 4.1
$ws.Range("I42").Value = 4.8
$ws.Range("J42").Value = 3.2
$ws.Range("K42").Value = 3.65
$ws.Range("L42").Value = 1.43
$ws.Range("M42").Value = 1.1
$ws.Range("N42").Value = 2.98
$ws.Range("O42").Value = 1.42
$ws.Range("P42").Value = 1.67
$ws.Range("Q42").Value = 2.24
$ws.Range("R42").Value = 1.25
$ws.Range("S42").Value = 4.2
$ws.Range("T42").Value = 1.96
$ws.Range("U42").Value = 1.87
$ws.Range("V42").Value = 1.26
$ws.Range("W42").Value = 1.84
$ws.Range("X42").Value = 11
$ws.Range("Y42").Value = 13.5
$ws.Range("Z42").Value = 34
$ws.Range("AB42").Value = 7.8
$ws.Range("AD42").Value = 19
$ws.Range("AE42").Value = 1000
$ws.Range("AF42").Value = 12.5
$ws.Range("AH42").Value = 22
$ws.Range("AJ42").Value = 1000
$ws.Range("AL42").Value = 1000
$ws.Range("AN42").Value = 21
$ws.Range("U44").Value = 1.95
$ws.Range("I45").Value = 4.2
$ws.Range("W45").Value = 1.73
$ws.Range("F46").Value = 1.63
$ws.Range("I46").Value = 8
$ws.Range("P46").Value = 1.77
$ws.Range("W46").Value = 2.36
$ws.Range("AF46").Value = 13
$ws.Range("AK46").Value = 29
$ws.Range("F47").Value = 2.96
$ws.Range("G47").Value = 3.3
$ws.Range("H47").Value = 2.32
$ws.Range("I47").Value = 2.5
$ws.Range("J47").Value = 3.6
$ws.Range("K47").Value = 4
$ws.Range("N47").Value = 1.1
$ws.Range("P47").Value = 2.14
$ws.Range("Q47").Value = 1.76
$ws.Range("R47").Value = 1.44
$ws.Range("S47").Value = 2.92
$ws.Range("T47").Value = 1.04
$ws.Range("U47").Value = 2.18
$ws.Range("V47").Value = 1.66
$ws.Range("X47").Value = 22
$ws.Range("Y47").Value = 15
$ws.Range("Z47").Value = 21
$ws.Range("AA47").Value = 40
$ws.Range("AB47").Value = 17.5
$ws.Range("AC47").Value = 10.5
$ws.Range("AD47").Value = 14
$ws.Range("AE47").Value = 30
$ws.Range("AF47").Value = 28
$ws.Range("AG47").Value = 16.5
$ws.Range("AH47").Value = 19.5
$ws.Range("AI47").Value = 42
$ws.Range("AJ47").Value = 65
$ws.Range("AK47").Value = 40
$ws.Range("AL47").Value = 48
$ws.Range("AN47").Value = 32
$ws.Range("AO47").Value = 20
$ws.Range("I48").Value = 2.26
$ws.Range("W48").Value = 1.36
$ws.Range("H49").Value = 3.45
$ws.Range("M49").Value = 1.07
$ws.Range("N49").Value = 3.45
$ws.Range("R49").Value = 1.33
$ws.Range("S49").Value = 3.55
$ws.Range("X49").Value = 16.5
$ws.Range("Y49").Value = 13.5
$ws.Range("AA49").Value = 70
$ws.Range("AB49").Value = 9.6
$ws.Range("AC49").Value = 8
$ws.Range("AD49").Value = 980
$ws.Range("AF49").Value = 980
$ws.Range("AG49").Value = 11.5
$ws.Range("AH49").Value = 980
$ws.Range("AI49").Value = 55
$ws.Range("AM49").Value = 130
$ws.Range("AN49").Value = 24
